$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# "Generate Report for handoff"
#
# A new handoff was produced for both the zh-cn and de-de locales:
#   - Status moves from "Handoff transform failed" to "Ready for handoff"
#   - A link to the freshly generated *.xlf handoff file is recorded
#   - The "Latest Handoff Datetime" is stamped
#   - The "Handoff Reason" flips from "Ignored" to "Include"
# -----------------------------------------------------------------------

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/bef722aedc1bf0aa2e82857d253a2e3a8bb11a7c"
$sourceId = "74618e38-2105-4f16-92d0-0b9f28d59713"
$handoffHash = "9c196dd0d46b3b0d7defaec533a98d79aed9e093"

$mdDisplay = "$sourceId.md"
$mdUrl = "$repoBase/e2e/$mdDisplay"
$configDisplay = ".localization-config"
$configUrl = "$repoBase/$configDisplay"

$locales = @(
    @{ Name = "zh-cn"; Datetime = "2016-01-13 02:39:15" },
    @{ Name = "de-de"; Datetime = "2016-01-13 02:39:35" }
)

# The Overview sheet mirrors each locale's status in column B (zh-cn) / C
# (de-de) - keep it in sync with the per-locale "Ready for handoff" status.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

foreach ($locale in $locales) {
    $sheetName = $locale.Name
    $handoffDatetime = $locale.Datetime

    $ws = $wb.Worksheets.Item($sheetName)

    $xlfName = "$sourceId.$handoffHash.$sheetName.xlf"
    $xlfUrl = "$repoBase/e2e/$xlfName"

    # B2: Status -> "Ready for handoff"
    $ws.Range("B2").Value = "Ready for handoff"

    # C2: new "Latest Handoff File" link to the generated xlf
    $ws.Range("C2").Value = $xlfName

    # D2: "Latest Handoff Datetime" gets stamped (keep the datetime display format)
    $ws.Range("D2").Value = $handoffDatetime
    $ws.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    # H2: "Handoff Reason" -> "Include"
    $ws.Range("H2").Value = "Include"

    # Rebuild the hyperlinks collection in display order (A2, C2, A3) so
    # relationship ids line up the way Excel would naturally assign them.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, "", "", $mdDisplay)
    $ws.Hyperlinks.Add($ws.Range("C2"), $xlfUrl, "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, "", "", $configDisplay)
}
